$d = $word.ActiveDocument

# --- Edit 1: the "Type" cell in the command-line-interface table had its
# example text split across several runs (" " + proofErr-wrapped "arg" +
# " -f filter"); collapse that run-fragmented text into a single run,
# dropping the spell-check proofErr markers, while leaving the preceding
# "Type" run (and its run formatting) untouched. ---
$d.Content.Find.Execute(" arg -f filter", $false, $false, $false, $false, `
    $false, $true, 1, $false, " arg -f filter", 2) | Out-Null

# --- Edit 2: append a new "<diagnostics-general>" Heading3 section (with its
# Symptom / Probable cause table, and a trailing empty Heading3 paragraph)
# right before the document's existing trailing empty paragraph, which must
# itself be left in place and unmodified. ---
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)
$insertRange = $lastPara.Range
$insertRange.Collapse(1)

# Split off a fresh, empty paragraph immediately before the document's real
# trailing paragraph, and use *that* new paragraph as the InsertXML anchor.
# InsertXML on a collapsed range splices new block content in at that point,
# but when the anchor point sits right at the end of the body (immediately
# before </w:body>), a trailing </w:tbl> in the injected XML can absorb the
# pre-existing final paragraph. Inserting into a dedicated spacer paragraph
# avoids disturbing the original trailing paragraph.
$insertRange.InsertParagraphBefore()
$anchorPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$anchorRange = $anchorPara.Range
$anchorRange.Collapse(1)

$xmlPayload = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>&lt;diagnostics-general&gt;</w:t></w:r><w:r><w:br/></w:r></w:p><w:tbl><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblpPr w:leftFromText="180" w:rightFromText="180" w:vertAnchor="text" w:tblpY="1"/><w:tblOverlap w:val="never"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="3145"/><w:gridCol w:w="3330"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="3145" w:type="dxa"/><w:shd w:val="clear" w:color="auto" w:fill="C5E0B3" w:themeFill="accent6" w:themeFillTint="66"/></w:tcPr><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="990"/></w:tabs><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Symptom</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3330" w:type="dxa"/><w:shd w:val="clear" w:color="auto" w:fill="C5E0B3" w:themeFill="accent6" w:themeFillTint="66"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Probable cause</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3145" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="HTMLPreformatted"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Clunking noise on bumps only</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3330" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="HTMLPreformatted"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Checking the struts</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3145" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="HTMLPreformatted"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Continuous clunking noise</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3330" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="HTMLPreformatted"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Checking the ball joints</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3145" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="HTMLPreformatted"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Ticks when in neutral</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3330" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="HTMLPreformatted"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Checking the exhaust</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3145" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="HTMLPreformatted"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Ticks only in reverse</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3330" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="HTMLPreformatted"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Checking the brakes</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3145" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="HTMLPreformatted"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Ticks in turns and curves</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3330" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="HTMLPreformatted"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Checking the CV joints</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3145" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="HTMLPreformatted"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>T</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">icks </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">only </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>when cold</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3330" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="HTMLPreformatted"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Checking the catalytic converter</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3145" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="HTMLPreformatted"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Ticks only at slow speed</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3330" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="HTMLPreformatted"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Checking the wheels</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$anchorRange.InsertXML($xmlPayload)
